# Add the new quotation row (2025-10-31) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$prevRow = 56
$newRow = 57

# Column A: date value (stored as Excel serial date). Copy the number
# format/style used by the preceding rows in column A so the new cell is
# formatted the same way as the rest of the date column.
$ws.Cells.Item($prevRow, 1).Copy() | Out-Null
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item($newRow, 1).Value = 45961

# Columns B:E: quotation values, stored as plain text strings (comma as the
# decimal separator), matching the existing rows in the sheet.
$ws.Cells.Item($newRow, 2).Value = "22,0341"
$ws.Cells.Item($newRow, 3).Value = "16,1343"
$ws.Cells.Item($newRow, 4).Value = "15,5326"
$ws.Cells.Item($newRow, 5).Value = "15,5326"
